# Edit language-05 diagram on slide 5: rename the four "Q5" labelled ovals
# (Oval 96, Oval 97, Oval 98, Oval 99) to Q7, Q9, Q6, Q8 respectively.
# The first "Q5" oval (Oval 74) is left unchanged.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)

$s.Shapes.Item("Oval 96").TextFrame.TextRange.Text = "Q7"
$s.Shapes.Item("Oval 97").TextFrame.TextRange.Text = "Q9"
$s.Shapes.Item("Oval 98").TextFrame.TextRange.Text = "Q6"
$s.Shapes.Item("Oval 99").TextFrame.TextRange.Text = "Q8"
